$d = $word.ActiveDocument
$nl = [char]11

# --- Paragraph 6: "Objetivos" body text ---
# was: Propiciar aos alunos... (moves down to Docente list, paragraph 8)
# now: the old "Programa resumido" summary sentence
$d.Paragraphs(6).Range.Text = "1. Princípios da corrosão. 2. Cinética da corrosão. 3. Formas de corrosão. 4. Proteção contra a corrosão. 5. Oxidação em temperaturas elevadas."

# --- Paragraph 8: "Docente(s) Responsável(eis)" list bullet ---
# was: 5817344 - Livia Melo Carneiro (moves down to Bibliografia body, paragraph 16)
# now: the old "Objetivos" paragraph text
$d.Paragraphs(8).Range.Text = "Propiciar aos alunos os conhecimentos básicos de corrosão, nos aspectos termodinâmicos e cinéticos, e descrever as principais formas de ataque e as técnicas de proteção contra a corrosão e a oxidação metálica."

# --- Paragraph 10: "Programa resumido" body text ---
# was: the summary sentence (moves up into paragraph 6)
# now: the long itemized "Programa" content (moves up from paragraph 12)
$d.Paragraphs(10).Range.Text = "1. Princípios da corrosão: Reações de oxi-redução. Potenciais de eletrodo - Sistema redox em estado de equilíbrio - Diagrama de Pourbaix" + $nl + "2. Cinética da corrosão: - Sistema redox em estado de não equilíbrio - Teoria do potencial misto  Passivação." + $nl + "3. Formas de corrosão: - Corrosão galvânica - Corrosão por pites e frestas - Corrosão intergranular - Corrosão sob tensão - Danos causados pelo hidrogênio." + $nl + "4. Proteção contra a corrosão: - Proteção catódica e anódica - Inibidores  Revestimentos." + $nl + "5. Oxidação em temperaturas elevadas - Fundamentos termodinâmicos - Mecanismos de transporte - Velocidade de oxidação - Oxidação de metais puros - Oxidação de ligas."

# --- Paragraph 12: "Programa" body text ---
# was: the long itemized content (moves up into paragraph 10)
# now: the "Método" evaluation sentence (moves up from paragraph 14)
$d.Paragraphs(12).Range.Text = "O aluno será avaliado através de duas provas escritas P1 e P2."

# --- Paragraph 14: "Avaliação" list bullet (bold labels + values) ---
# Locate each bold label dynamically, then overwrite only the value run that
# follows it, leaving "Método:"/"Critério:"/"Norma de recuperação:" untouched.
$para = $d.Paragraphs(14).Range
$paraEnd = $para.End

$f1 = $d.Range($para.Start, $paraEnd)
$f1.Find.Execute("Método: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lbl1End = $f1.End

$f2 = $d.Range($lbl1End, $paraEnd)
$f2.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lbl2Start = $f2.Start
$lbl2End = $f2.End

$f3 = $d.Range($lbl2End, $paraEnd)
$f3.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lbl3Start = $f3.Start
$lbl3End = $f3.End

$val1 = $d.Range($lbl1End, $lbl2Start - 1)
$val2 = $d.Range($lbl2End, $lbl3Start - 1)
$val3 = $d.Range($lbl3End, $paraEnd)

# Assign back-to-front so already-computed offsets stay valid.
# was: Método->"O aluno será avaliado..."; Critério->"A nota final NF será...";
#      Norma->"Prova escrita sobre toda matéria..."
# now: Método->"A nota final NF será..."; Critério->"Prova escrita sobre toda matéria...";
#      Norma->the old Bibliografia list (moves up from paragraph 16)
$val3.Text = "V.GENTIL, Corrosão, Ed. Guanabara Dois, 1982" + $nl + "L.V. RAMANATHAN, Corrosão e seu controle, Ed. Hermes" + $nl + "L.L. SHREIR, Corrosion, Newnes Butterworths, 2 vol., 1976" + $nl + "N. BIRKS and G.H.MEIER, Introduction to High Temperature Oxidation of Metals, Edward Arnold, 1983"
$val2.Text = "Prova escrita sobre toda matéria. A média final MF será a média da nota final NF e da nota obtida na recuperação NR: MF = (NF + NR)/2 ."
$val1.Text = "A nota final NF será (P1 + P2)/2 ."

# --- Paragraph 16: "Bibliografia" body text ---
# was: the bibliography list (moves up into paragraph 14's Norma value)
# now: the docente name (moves up from paragraph 8)
$d.Paragraphs(16).Range.Text = "5817344 - Livia Melo Carneiro"
